# Scheduled-runner update: refresh cached market-board pricing figures
# (currentAveragePrice*, LevePrice*, LeveProfit*) across the per-job
# leve-profit sheets. Values only; no structural changes.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H33").Value = 404
$ws.Range("I33").Value = 220.15384
$ws.Range("J33").Value = 1200.6666
$ws.Range("K33").Value = 220.15384
$ws.Range("L33").Value = 1200.6666
$ws.Range("M33").Value = 8.846159999999998
$ws.Range("N33").Value = -1658.6666

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H112").Value = 54104.184
$ws.Range("I112").Value = 1286.625
$ws.Range("J112").Value = 84285.64
$ws.Range("K112").Value = 3859.875
$ws.Range("L112").Value = 252856.92
$ws.Range("M112").Value = -2751.875
$ws.Range("N112").Value = -255072.92

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H137").Value = 328832.06
$ws.Range("I137").Value = 689289.3
$ws.Range("J137").Value = 10781.529
$ws.Range("K137").Value = 2067867.9
$ws.Range("L137").Value = 32344.587
$ws.Range("M137").Value = -2065317.9
$ws.Range("N137").Value = -37444.587

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H138").Value = 5559.5894
$ws.Range("I138").Value = 1510.1666
$ws.Range("J138").Value = 6663.977
$ws.Range("K138").Value = 4530.4998
$ws.Range("L138").Value = 19991.931
$ws.Range("M138").Value = 609.5002000000004
$ws.Range("N138").Value = -30271.931

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H141").Value = 5906.095
$ws.Range("I141").Value = 5968.5386
$ws.Range("J141").Value = 5804.625
$ws.Range("K141").Value = 17905.6158
$ws.Range("L141").Value = 17413.875
$ws.Range("M141").Value = -12725.6158
$ws.Range("N141").Value = -27773.875

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H122").Value = 440883.9
$ws.Range("I122").Value = 2520.9412
$ws.Range("J122").Value = 937695.25
$ws.Range("K122").Value = 7562.823600000001
$ws.Range("L122").Value = 2813085.75
$ws.Range("M122").Value = -5112.823600000001
$ws.Range("N122").Value = -2817985.75

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H132").Value = 2751.4138
$ws.Range("I132").Value = 2080.5715
$ws.Range("J132").Value = 4512.375
$ws.Range("K132").Value = 6241.7145
$ws.Range("L132").Value = 13537.125
$ws.Range("M132").Value = -3711.7145
$ws.Range("N132").Value = -18597.125

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H105").Value = 33756.742
$ws.Range("I105").Value = 38205.48
$ws.Range("J105").Value = 3727.75
$ws.Range("K105").Value = 38205.48
$ws.Range("L105").Value = 3727.75
$ws.Range("M105").Value = -36458.48
$ws.Range("N105").Value = -7221.75

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 2721
$ws.Range("I31").Value = 1893.8125
$ws.Range("J31").Value = 3147.9355
$ws.Range("K31").Value = 1893.8125
$ws.Range("L31").Value = 3147.9355
$ws.Range("M31").Value = -1598.8125
$ws.Range("N31").Value = -3737.9355

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H34").Value = 2721
$ws.Range("I34").Value = 1893.8125
$ws.Range("J34").Value = 3147.9355
$ws.Range("K34").Value = 1893.8125
$ws.Range("L34").Value = 3147.9355
$ws.Range("M34").Value = -1691.8125
$ws.Range("N34").Value = -3551.9355

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H68").Value = 45763.5
$ws.Range("J68").Value = 45763.5
$ws.Range("L68").Value = 45763.5
$ws.Range("N68").Value = -47261.5

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H71").Value = 45763.5
$ws.Range("J71").Value = 45763.5
$ws.Range("L71").Value = 137290.5
$ws.Range("N71").Value = -144778.5

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H74").Value = 45833
$ws.Range("J74").Value = 39500
$ws.Range("L74").Value = 39500
$ws.Range("N74").Value = -41248

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H77").Value = 45833
$ws.Range("J77").Value = 39500
$ws.Range("L77").Value = 118500
$ws.Range("N77").Value = -127236

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H132").Value = 23835258
$ws.Range("I132").Value = 25007020
$ws.Range("J132").Value = 399999.5
$ws.Range("K132").Value = 75021060
$ws.Range("L132").Value = 1199998.5
$ws.Range("M132").Value = -75018530
$ws.Range("N132").Value = -1205058.5

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H137").Value = 20000
$ws.Range("J137").Value = 0
$ws.Range("L137").Value = 0
$ws.Range("N137").ClearContents()

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H138").Value = 45000
$ws.Range("J138").Value = 65000
$ws.Range("L138").Value = 65000
$ws.Range("N138").Value = -75280

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H140").Value = 58488
$ws.Range("J140").Value = 86976
$ws.Range("L140").Value = 86976
$ws.Range("N140").Value = -97336

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H68").Value = 10491.579
$ws.Range("J68").Value = 15859.546
$ws.Range("L68").Value = 47578.638
$ws.Range("N68").Value = -49200.638

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H71").Value = 10491.579
$ws.Range("J71").Value = 15859.546
$ws.Range("L71").Value = 142735.914
$ws.Range("N71").Value = -150847.914

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H132").Value = 9288000
$ws.Range("I132").Value = 668.375
$ws.Range("J132").Value = 16717865
$ws.Range("K132").Value = 6015.375
$ws.Range("L132").Value = 150460785
$ws.Range("M132").Value = -3485.375
$ws.Range("N132").Value = -150465845

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H132").Value = 2917.0425
$ws.Range("I132").Value = 2551.1082
$ws.Range("J132").Value = 4271
$ws.Range("K132").Value = 7653.3246
$ws.Range("L132").Value = 12813
$ws.Range("M132").Value = -5123.3246
$ws.Range("N132").Value = -17873

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H135").Value = 0
$ws.Range("J135").Value = 0
$ws.Range("L135").Value = 0
$ws.Range("N135").ClearContents()

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H132").Value = 4202065
$ws.Range("I132").Value = 4964895
$ws.Range("J132").Value = 6500
$ws.Range("K132").Value = 14894685
$ws.Range("L132").Value = 19500
$ws.Range("M132").Value = -14892155
$ws.Range("N132").Value = -24560

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H81").Value = 56398.2
$ws.Range("I81").Value = 56398.2
$ws.Range("J81").Value = 0
$ws.Range("K81").Value = 112796.4
$ws.Range("L81").Value = 0
$ws.Range("M81").Value = -111735.4
$ws.Range("N81").ClearContents()

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H84").Value = 56398.2
$ws.Range("I84").Value = 56398.2
$ws.Range("J84").Value = 0
$ws.Range("K84").Value = 563982
$ws.Range("L84").Value = 0
$ws.Range("M84").Value = -558678
$ws.Range("N84").ClearContents()

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value = 24554.2
$ws.Range("I132").Value = 27468.3
$ws.Range("J132").Value = 12897.8
$ws.Range("K132").Value = 82404.9
$ws.Range("L132").Value = 38693.39999999999
$ws.Range("M132").Value = -79874.9
$ws.Range("N132").Value = -43753.39999999999
